# Rework the simple task list into a 3-column sheet that also tracks
# a related file attachment, and drop the now-unused rows of sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Tasks"
$ws.Range("B1").Value = "Due Date"
$ws.Range("C1").Value = "Related Materials (file name)"

# Single remaining task row
$ws.Range("A2").Value = "[ ] CSM"
$ws.Range("B2").Value = "Friday"
$ws.Range("C2").Value = "week13.pdf"

# The rest of the old sample tasks/dates are removed, but keep the bold
# task-row formatting (style) on column A for the now-empty rows.
$ws.Range("A3:A6").ClearContents()
$ws.Range("B3:B6").ClearContents()

# New "Related Materials" column needs to be wide enough to show the
# file name / header text.
$ws.Columns("C").ColumnWidth = 23.83

# Leave the selection where it was left after adding the new column.
[void]$ws.Range("C3").Select()
